$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 47, shifting existing rows 47-103 down to 48-104.
$ws.Rows.Item(47).Insert(-4121)

# Populate the newly inserted row 47 with its data.
$ws.Range("A47").Value = 7
$ws.Range("B47").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C47").Value = "Ñuble"
$ws.Range("D47").Value = 44944
$ws.Range("D47").NumberFormat = $ws.Range("D48").NumberFormat
$ws.Range("E47").Value = 16
$ws.Range("F47").Value = 100112022
$ws.Range("G47").Value = "Arveja Verde"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 50
$ws.Range("K47").Value = 22000
$ws.Range("L47").Value = 22000
$ws.Range("M47").Value = 22000
$ws.Range("N47").Value = "$/saco 25 kilos"
$ws.Range("O47").Value = "Región de Ñuble"
$ws.Range("P47").Value = 880
$ws.Range("Q47").Value = 25
$ws.Range("R47").Value = "Hortaliza"
